$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value that LOOKS like a date ("DD-MMM-YY") into a cell as
# plain text, without letting Excel's autoconvert turn it into a date serial.
# We type it with a leading apostrophe (forces text / quotePrefix), then fix
# the cell's formatting back up by pasting the format from a known-good
# "plain text" cell in the same row (so the final style index matches the
# rest of the row instead of picking up a stray quotePrefix style).
# ---------------------------------------------------------------------------
function Set-DateLikeText($addr, $val, $fmtSourceAddr) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($fmtSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Grow the table from 4 data rows to 6: clone row 4's formatting down into
# the two new rows (5 and 6) before filling in values, so every new cell
# picks up the correct style (s=2 body / s=3 "LOW THREAT" impact cell).
# ---------------------------------------------------------------------------
$ws.Range("A4:K4").Copy()
$ws.Range("A5:K6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 2: 26-FEB-26 / SM-322 / Nile Air NP-120 (unchanged) — only the fare
# figures move.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 650
$ws.Range("E2").Value = 728
$ws.Range("F2").Value = -78

# ---------------------------------------------------------------------------
# Row 3: becomes 05-MAR-26 / SM-322 / Nile Air NP-120, new fares, and the
# impact flips from HIGH THREAT ALERT (red) to LOW THREAT (green).
# ---------------------------------------------------------------------------
Set-DateLikeText "A3" "05-MAR-26" "B3"
$ws.Range("B3").Value = "SM-322"
$ws.Range("C3").Value = "Nile Air NP-120"
$ws.Range("D3").Value = 745
$ws.Range("E3").Value = 848
$ws.Range("F3").Value = -103
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("J2").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 4: becomes 11-MAR-26 / SM-970 / Nile Air NP-110, new fares (impact
# stays LOW THREAT).
# ---------------------------------------------------------------------------
Set-DateLikeText "A4" "11-MAR-26" "B4"
$ws.Range("B4").Value = "SM-970"
$ws.Range("C4").Value = "Nile Air NP-110"
$ws.Range("D4").Value = 745
$ws.Range("E4").Value = 788
$ws.Range("F4").Value = -43

# ---------------------------------------------------------------------------
# Row 5 (new): 12-MAR-26 / SM-322 / Nile Air NP-120.
# ---------------------------------------------------------------------------
Set-DateLikeText "A5" "12-MAR-26" "B5"
$ws.Range("B5").Value = "SM-322"
$ws.Range("C5").Value = "Nile Air NP-120"
$ws.Range("D5").Value = 845
$ws.Range("E5").Value = 848
$ws.Range("F5").Value = -3
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

# ---------------------------------------------------------------------------
# Row 6 (new): 26-MAR-26 / SM-322 / Nile Air NP-120.
# ---------------------------------------------------------------------------
Set-DateLikeText "A6" "26-MAR-26" "B6"
$ws.Range("B6").Value = "SM-322"
$ws.Range("C6").Value = "Nile Air NP-120"
$ws.Range("D6").Value = 745
$ws.Range("E6").Value = 788
$ws.Range("F6").Value = -43
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

# ---------------------------------------------------------------------------
# Column J (10th column) narrows from 30 to 12 characters wide now that the
# long "HIGH THREAT ALERT - NEED ACTION" text is gone. Excel's ColumnWidth
# setter stores width + 5/6 in the saved file, so back that offset out.
# ---------------------------------------------------------------------------
$ws.Range("J1").ColumnWidth = 12 - 5/6
